$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-15 Friday" "2024-03-16 Saturday"
Replace-Text "391×3=1173" "977×2=1954"
Replace-Text "816×7=5712" "309×4=1236"
Replace-Text "874×7=6118" "536×7=3752"
Replace-Text "936×5=4680" "359×8=2872"
Replace-Text "177×8=1416" "605×7=4235"
Replace-Text "173×2=346" "606×2=1212"
Replace-Text "744×7=5208" "548×4=2192"
Replace-Text "770×7=5390" "130×9=1170"
Replace-Text "971×8=7768" "449×2=898"
Replace-Text "390×3=1170" "239×9=2151"
Replace-Text "551×6=3306" "138×5=690"
Replace-Text "917×4=3668" "143×6=858"
Replace-Text "248×8=1984" "577×8=4616"
Replace-Text "248×7=1736" "341×9=3069"
Replace-Text "826×3=2478" "903×6=5418"
Replace-Text "486×7=3402" "928×4=3712"
Replace-Text "910×4=3640" "610×4=2440"
Replace-Text "645×3=1935" "313×8=2504"
Replace-Text "627×8=5016" "468×2=936"
Replace-Text "394×3=1182" "592×8=4736"
Replace-Text "307×4=1228" "661×9=5949"
Replace-Text "668×8=5344" "955×9=8595"
Replace-Text "365×4=1460" "583×7=4081"
Replace-Text "497×9=4473" "466×8=3728"
Replace-Text "925×8=7400" "297×2=594"
